$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Beginner" tree-pose button label/value to the new
# "Label:ButtonName,Label:ButtonName" style string used for dynamically
# creating buttons on the form.
$ws.Range("B1").Value = "Tree Pose:TreePose,Cat Pose:CatPose"

# Widen column B to fit the new button-config string (~39.71 chars wide).
$ws.Columns.Item(2).ColumnWidth = 38.75

# Update the saved selection / scroll position for the sheet.
$ws.Range("E5").Select()
